$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Random Forest, Logistic Regression, Support Vector Machine." ->
#    "Random Forest, Logistic Regression and Support Vector Machine."
#    (only the occurrence that is followed by "The models will classify")
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Regression, Support Vector Machine. The models will classify",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Regression and Support Vector Machine. The models will classify", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Add a first-line indent (720 twips = 36 pt) to the paragraph that
#    begins "The dataset is taken from the Kaggle repository..."
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "*The dataset is taken from the Kaggle repository*includes boolean*") {
        $p.Range.ParagraphFormat.FirstLineIndent = 36
        break
    }
}

# ------------------------------------------------------------------
# 3. "The dataset includes boolean, float, int, and string types." ->
#    "The dataset contains numeric values and Boolean values. There
#    are missing values in the dataset."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "The dataset includes boolean, float, int, and string types.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The dataset contains numeric values and Boolean values. There are missing values in the dataset.", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "...TBG, referral_source." -> "...TBG, and referral_source."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "TBG, referral_source.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "TBG, and referral_source.", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Remove the first-line indent from the paragraph that begins
#    "The class labels include letters from A to T..."
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "The class labels include letters from A to T*") {
        $p.Range.ParagraphFormat.FirstLineIndent = 0
        break
    }
}

# ------------------------------------------------------------------
# 6. Drop the trailing space after "...different thyroid conditions."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "different thyroid conditions. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "different thyroid conditions.", 2) | Out-Null
